$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.797.03'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +3.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.677.36'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.33%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9980'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.65'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9996'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4631'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.27%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2604'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06147'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.674.10'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.13%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06996'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.36%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.95'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.26%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.366'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5785'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '75.47'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9996'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9985'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.20%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.783.86'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +3.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006709'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.44'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.887.59'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.470'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +3.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.672'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.229'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.23'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.87%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.99'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.71%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.716'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +5.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '104.76'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.950'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.51%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07694'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.613'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04362'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +2.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.598'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6082'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.76%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9534'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9321'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '109.33'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +10.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.453'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9986'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.867'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01456'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.35%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.049'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +8.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3728'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1120'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.05%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.97%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.158'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.32%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '31.27'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +9.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.623'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +7.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.214'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.85%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.02%  '
